# Repull data: update the "dSF" (column F) values for a set of rows to
# reflect the newly pulled data / recalculated mean.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 0
    12 = 6
    13 = 2
    28 = 1
    35 = 0
    36 = 0
    37 = -1
    39 = -6
    49 = -1
    57 = -2
    66 = 0
    67 = 2
    76 = -7
    77 = -2
    79 = -1
    82 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
